$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $p.Range
$start = $r.Start
$r.InsertAfter("SOME OTHER THING")
$newRange = $d.Range($start, $start + 16)
$newRange.LanguageID = "en-US"
